$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows that are removed from the left block (A:H) and the trailing row (J:Q)
$ws.Range("A10:H12").Clear()
$ws.Range("J35:Q35").Clear()

# Update changed cell values
$ws.Range("A3").Value = "crude"
$ws.Range("B3").Value = 0.7647058823529411
$ws.Range("C3").Value = 26
$ws.Range("D3").Value = 26
$ws.Range("H3").Value = 8
$ws.Range("K3").Value = 0.9696969696969697
$ws.Range("L3").Value = 32
$ws.Range("M3").Value = 32
$ws.Range("Q3").Value = 1
$ws.Range("A4").Value = "died"
$ws.Range("B4").Value = 0.7307692307692307
$ws.Range("C4").Value = 19
$ws.Range("D4").Value = 19
$ws.Range("H4").Value = 7
$ws.Range("J4").Value = "best"
$ws.Range("K4").Value = 0.9322033898305084
$ws.Range("L4").Value = 55
$ws.Range("M4").Value = 55
$ws.Range("Q4").Value = 4
$ws.Range("A5").Value = "fraud"
$ws.Range("B5").Value = 0.6111111111111112
$ws.Range("C5").Value = 22
$ws.Range("D5").Value = 22
$ws.Range("H5").Value = 14
$ws.Range("J5").Value = "love"
$ws.Range("K5").Value = 0.8913043478260869
$ws.Range("L5").Value = 41
$ws.Range("M5").Value = 41
$ws.Range("Q5").Value = 5
$ws.Range("B6").Value = 0.5753424657534246
$ws.Range("C6").Value = 168
$ws.Range("D6").Value = 168
$ws.Range("H6").Value = 124
$ws.Range("J6").Value = "happy"
$ws.Range("K6").Value = 0.8846153846153846
$ws.Range("L6").Value = 23
$ws.Range("M6").Value = 23
$ws.Range("Q6").Value = 3
$ws.Range("A7").Value = "panic"
$ws.Range("B7").Value = 0.1782945736434109
$ws.Range("C7").Value = 92
$ws.Range("D7").Value = 92
$ws.Range("H7").Value = 424
$ws.Range("K7").Value = 0.8392857142857143
$ws.Range("L7").Value = 94
$ws.Range("M7").Value = 94
$ws.Range("Q7").Value = 18
$ws.Range("A8").Value = "sc"
$ws.Range("B8").Value = 0.1587301587301587
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = 30
$ws.Range("H8").Value = 159
$ws.Range("J8").Value = "won"
$ws.Range("K8").Value = 0.8205128205128205
$ws.Range("L8").Value = 32
$ws.Range("M8").Value = 32
$ws.Range("Q8").Value = 7
$ws.Range("A9").Value = "no"
$ws.Range("B9").Value = 0.05555555555555555
$ws.Range("C9").Value = 20
$ws.Range("D9").Value = 20
$ws.Range("H9").Value = 340
$ws.Range("J9").Value = "thank"
$ws.Range("K9").Value = 0.8046875
$ws.Range("L9").Value = 103
$ws.Range("M9").Value = 103
$ws.Range("Q9").Value = 25
$ws.Range("J10").Value = "thanks"
$ws.Range("K10").Value = 0.7926829268292683
$ws.Range("L10").Value = 65
$ws.Range("M10").Value = 65
$ws.Range("Q10").Value = 17
$ws.Range("J11").Value = "free"
$ws.Range("K11").Value = 0.7583333333333333
$ws.Range("L11").Value = 91
$ws.Range("M11").Value = 91
$ws.Range("Q11").Value = 29
$ws.Range("J12").Value = "special"
$ws.Range("K12").Value = 0.75
$ws.Range("L12").Value = 27
$ws.Range("M12").Value = 27
$ws.Range("Q12").Value = 9
$ws.Range("J13").Value = "positive"
$ws.Range("K13").Value = 0.7413793103448276
$ws.Range("Q13").Value = 15
$ws.Range("J14").Value = "nice"
$ws.Range("K14").Value = 0.7407407407407407
$ws.Range("L14").Value = 20
$ws.Range("M14").Value = 20
$ws.Range("Q14").Value = 7
$ws.Range("J15").Value = "confidence"
$ws.Range("K15").Value = 0.7222222222222222
$ws.Range("L15").Value = 26
$ws.Range("M15").Value = 26
$ws.Range("Q15").Value = 10
$ws.Range("J16").Value = "safety"
$ws.Range("K16").Value = 0.7058823529411765
$ws.Range("L16").Value = 36
$ws.Range("M16").Value = 36
$ws.Range("Q16").Value = 15
$ws.Range("K17").Value = 0.704225352112676
$ws.Range("L17").Value = 100
$ws.Range("M17").Value = 100
$ws.Range("Q17").Value = 42
$ws.Range("J18").Value = "good"
$ws.Range("K18").Value = 0.65625
$ws.Range("L18").Value = 105
$ws.Range("M18").Value = 105
$ws.Range("Q18").Value = 55
$ws.Range("J19").Value = "support"
$ws.Range("K19").Value = 0.6415094339622641
$ws.Range("L19").Value = 68
$ws.Range("M19").Value = 68
$ws.Range("Q19").Value = 38
$ws.Range("J20").Value = "relief"
$ws.Range("K20").Value = 0.64
$ws.Range("L20").Value = 32
$ws.Range("M20").Value = 32
$ws.Range("Q20").Value = 18
$ws.Range("K21").Value = 0.6349206349206349
$ws.Range("L21").Value = 40
$ws.Range("M21").Value = 40
$ws.Range("Q21").Value = 23
$ws.Range("J22").Value = "well"
$ws.Range("K22").Value = 0.6170212765957447
$ws.Range("L22").Value = 58
$ws.Range("M22").Value = 58
$ws.Range("Q22").Value = 36
$ws.Range("J23").Value = "heroes"
$ws.Range("K23").Value = 0.5957446808510638
$ws.Range("L23").Value = 28
$ws.Range("M23").Value = 28
$ws.Range("Q23").Value = 19
$ws.Range("J24").Value = "fresh"
$ws.Range("K24").Value = 0.5416666666666666
$ws.Range("L24").Value = 26
$ws.Range("M24").Value = 26
$ws.Range("Q24").Value = 22
$ws.Range("K25").Value = 0.5195822454308094
$ws.Range("L25").Value = 199
$ws.Range("M25").Value = 199
$ws.Range("Q25").Value = 184
$ws.Range("K26").Value = 0.4352941176470588
$ws.Range("L26").Value = 148
$ws.Range("M26").Value = 148
$ws.Range("Q26").Value = 192
$ws.Range("J27").Value = "protect"
$ws.Range("K27").Value = 0.4246575342465753
$ws.Range("L27").Value = 31
$ws.Range("M27").Value = 31
$ws.Range("Q27").Value = 42
$ws.Range("J28").Value = "help"
$ws.Range("K28").Value = 0.4169491525423729
$ws.Range("L28").Value = 123
$ws.Range("M28").Value = 123
$ws.Range("Q28").Value = 172
$ws.Range("J29").Value = "care"
$ws.Range("K29").Value = 0.4157303370786517
$ws.Range("L29").Value = 37
$ws.Range("M29").Value = 37
$ws.Range("Q29").Value = 52
$ws.Range("J30").Value = "please"
$ws.Range("K30").Value = 0.3472803347280335
$ws.Range("L30").Value = 83
$ws.Range("M30").Value = 83
$ws.Range("Q30").Value = 156
$ws.Range("J31").Value = "increase"
$ws.Range("K31").Value = 0.3076923076923077
$ws.Range("L31").Value = 24
$ws.Range("M31").Value = 24
$ws.Range("Q31").Value = 54
$ws.Range("J32").Value = "and"
$ws.Range("K32").Value = 0.008614232209737827
$ws.Range("L32").Value = 23
$ws.Range("M32").Value = 26
$ws.Range("N32").Value = 0.88
$ws.Range("O32").Value = 0.12
$ws.Range("P32").Value = $true
$ws.Range("Q32").Value = 2647
$ws.Range("J33").Value = "."
$ws.Range("K33").Value = 0.005406487785342411
$ws.Range("L33").Value = 27
$ws.Range("M33").Value = 29
$ws.Range("N33").Value = 0.93
$ws.Range("O33").Value = 0.06999999999999995
$ws.Range("P33").Value = $true
$ws.Range("Q33").Value = 4967
$ws.Range("J34").Value = "to"
$ws.Range("K34").Value = 0.005314232902033271
$ws.Range("L34").Value = 23
$ws.Range("M34").Value = 23
$ws.Range("Q34").Value = 4305
